# Update the sentiment labels and their counts to reflect the new,
# more robust UiPath selector-based comment retrieval results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Very Negative"
$ws.Range("B1").Value = 5

$ws.Range("A2").Value = "Neutral"
$ws.Range("B2").Value = 67

$ws.Range("A3").Value = "Positive"
$ws.Range("B3").Value = 16

$ws.Range("A4").Value = "Very Positive"
$ws.Range("B4").Value = 5

$ws.Range("A5").Value = "Negative"
$ws.Range("B5").Value = 7
